$d = $word.ActiveDocument

function InsertRun($text, $bold) {
    $rng = $d.Range($global:pos, $global:pos)
    $rng.InsertAfter($text)
    if ($text.Length -gt 0) {
        $rng.Font.Bold = $bold
    }
    $global:pos = $rng.End
}

# Step 1: locate and delete the "License Information" heading paragraph entirely.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "License Information") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

# Step 2: locate the big "is based on" / credits paragraph and the following
# "This PDF version..." paragraph, then merge them by deleting the latter and
# rewriting the former's runs.
$bigPara = $null
$pdfPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.Contains("Translation Questions") -and $t.Contains("is based on")) {
        $bigPara = $p
    } elseif ($t.Contains("This PDF version is provided under the same license")) {
        $pdfPara = $p
    }
}

if ($pdfPara -ne $null) {
    $pdfPara.Range.Delete()
}

$full = $bigPara.Range
$clearRng = $d.Range($full.Start, $full.End - 1)
$clearRng.Text = ""
$global:pos = $clearRng.Start

InsertRun "" $false
InsertRun "unfoldingWord® Translation Questions" $true
InsertRun " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. " $false
InsertRun "unfoldingWord® Translation Questions" $false
InsertRun " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from " $false
InsertRun "unfoldingWord® Translation Questions" $false
InsertRun " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual" $false
InsertRun "" $false
